$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 (A=0) ---
$ws.Range("B2").Value = 1.901710291787398
$ws.Range("C2").Value = 4185.8
$ws.Range("D2").Value = 0.01349240180076153
$ws.Range("E2").Value = 29.6
$ws.Range("F2").Value = 204.4
$ws.Range("G2").Value = "MyDogaN"
$ws.Range("H2").Value = "DUO_SUPPORT"
$ws.Range("I2").Value = 0.09128055152368671
$ws.Range("J2").Value = 18.6
$ws.Range("K2").Value = 0.008167217339014521

# --- Update existing row 3 (A=1) ---
$ws.Range("B3").Value = 2.791646791513082
$ws.Range("C3").Value = 3890.6
$ws.Range("D3").Value = 0.01949175076877424
$ws.Range("E3").Value = 28.2
$ws.Range("F3").Value = 221.4
$ws.Range("G3").Value = "Mr Kayn"
$ws.Range("H3").Value = "DUO_SUPPORT"
$ws.Range("I3").Value = 0.1749684198889241
$ws.Range("J3").Value = 3.2
$ws.Range("K3").Value = 0.002463335941977546

# --- Insert three new rows (4,5,6) below, inheriting row-3 style, then fix borders ---
$ws.Rows("4:6").Insert()
$ws.Range("A4:A6").Borders.LineStyle = 1

# Row 4 (A=2) - all-zero / blank SOLO row
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("H4").Value = "SOLO"
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0

# Row 5 (A=3)
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 1.775319622012229
$ws.Range("C5").Value = 3193.8
$ws.Range("D5").Value = 0.01634241245136187
$ws.Range("E5").Value = 29.4
$ws.Range("F5").Value = 41
$ws.Range("G5").Value = "Portgas D " + [char]0x00C5 + "ce "
$ws.Range("H5").Value = "SOLO"
$ws.Range("I5").Value = 0.02279043913285158
$ws.Range("J5").Value = 4.4
$ws.Range("K5").Value = 0.002445803224013341

# Row 6 (A=4)
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 7.410958116892823
$ws.Range("C6").Value = 13175.2
$ws.Range("D6").Value = 0.05164043510720542
$ws.Range("E6").Value = 91.2
$ws.Range("F6").Value = 278.8
$ws.Range("G6").Value = "LS DUFFY"
$ws.Range("H6").Value = "SOLO"
$ws.Range("I6").Value = 0.1567269453783251
$ws.Range("J6").Value = 8.199999999999999
$ws.Range("K6").Value = 0.004627791059971394

Write-Output "edit applied"
